$d = $word.ActiveDocument

$d.Content.Find.Execute("161×4=", $true, $false, $false, $false, $false, $true, 1, $false, "241×4=", 2) | Out-Null
$d.Content.Find.Execute("643×4=", $true, $false, $false, $false, $false, $true, 1, $false, "318×2=", 2) | Out-Null
$d.Content.Find.Execute("658×7=", $true, $false, $false, $false, $false, $true, 1, $false, "225×5=", 2) | Out-Null
$d.Content.Find.Execute("570×7=", $true, $false, $false, $false, $false, $true, 1, $false, "451×2=", 2) | Out-Null
$d.Content.Find.Execute("156×6=", $true, $false, $false, $false, $false, $true, 1, $false, "822×6=", 2) | Out-Null
$d.Content.Find.Execute("779×9=", $true, $false, $false, $false, $false, $true, 1, $false, "113×2=", 2) | Out-Null
$d.Content.Find.Execute("655×6=", $true, $false, $false, $false, $false, $true, 1, $false, "712×2=", 2) | Out-Null
$d.Content.Find.Execute("859×5=", $true, $false, $false, $false, $false, $true, 1, $false, "255×5=", 2) | Out-Null
$d.Content.Find.Execute("905×2=", $true, $false, $false, $false, $false, $true, 1, $false, "470×6=", 2) | Out-Null
$d.Content.Find.Execute("744×6=", $true, $false, $false, $false, $false, $true, 1, $false, "571×7=", 2) | Out-Null
$d.Content.Find.Execute("299×8=", $true, $false, $false, $false, $false, $true, 1, $false, "517×6=", 2) | Out-Null
$d.Content.Find.Execute("933×2=", $true, $false, $false, $false, $false, $true, 1, $false, "773×5=", 2) | Out-Null
$d.Content.Find.Execute("249×6=", $true, $false, $false, $false, $false, $true, 1, $false, "141×5=", 2) | Out-Null
$d.Content.Find.Execute("665×5=", $true, $false, $false, $false, $false, $true, 1, $false, "745×4=", 2) | Out-Null
$d.Content.Find.Execute("440×7=", $true, $false, $false, $false, $false, $true, 1, $false, "390×5=", 2) | Out-Null
$d.Content.Find.Execute("515×7=", $true, $false, $false, $false, $false, $true, 1, $false, "307×3=", 2) | Out-Null
$d.Content.Find.Execute("486×7=", $true, $false, $false, $false, $false, $true, 1, $false, "343×2=", 2) | Out-Null
$d.Content.Find.Execute("975×2=", $true, $false, $false, $false, $false, $true, 1, $false, "954×2=", 2) | Out-Null
$d.Content.Find.Execute("545×5=", $true, $false, $false, $false, $false, $true, 1, $false, "465×5=", 2) | Out-Null
$d.Content.Find.Execute("233×6=", $true, $false, $false, $false, $false, $true, 1, $false, "883×7=", 2) | Out-Null
$d.Content.Find.Execute("528×4=", $true, $false, $false, $false, $false, $true, 1, $false, "643×4=", 2) | Out-Null
$d.Content.Find.Execute("776×2=", $true, $false, $false, $false, $false, $true, 1, $false, "991×6=", 2) | Out-Null
$d.Content.Find.Execute("938×4=", $true, $false, $false, $false, $false, $true, 1, $false, "829×7=", 2) | Out-Null
$d.Content.Find.Execute("122×8=", $true, $false, $false, $false, $false, $true, 1, $false, "475×8=", 2) | Out-Null
$d.Content.Find.Execute("114×7=", $true, $false, $false, $false, $false, $true, 1, $false, "828×8=", 2) | Out-Null
